$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: refresh existing values in place (ticker M.BA -> BABA.BA, new metrics). ---
# Column F (initial_date) is left untouched - it is unchanged by the edit.
$ws.Cells.Item(2, 1).Value = "BABA.BA"
$ws.Cells.Item(2, 2).Value = 20
$ws.Cells.Item(2, 3).Value = 7
$ws.Cells.Item(2, 4).Value = 1000000
$ws.Cells.Item(2, 5).Value = 1359390.288
$ws.Cells.Item(2, 7).Value = 45190
$ws.Cells.Item(2, 8).Value = 23
$ws.Cells.Item(2, 9).Value = 8
$ws.Cells.Item(2, 10).Value = 15
$ws.Cells.Item(2, 11).Value = 0.348
$ws.Cells.Item(2, 12).Value = 0.652
$ws.Cells.Item(2, 13).Value = -0.158
$ws.Cells.Item(2, 14).Value = 170013.286
$ws.Cells.Item(2, 15).Value = 0.208
$ws.Cells.Item(2, 16).Value = -94548.67200000001
$ws.Cells.Item(2, 17).Value = -0.101
$ws.Cells.Item(2, 18).Value = -145094.332
$ws.Cells.Item(2, 19).Value = 351006.288
$ws.Cells.Item(2, 20).Value = 0.3510062880000002
$ws.Cells.Item(2, 21).Value = 3
$ws.Cells.Item(2, 22).Value = 7
$ws.Cells.Item(2, 23).Value = 53.36363636363637

# --- New rows 3-5: additional tickers analyzed with the same strategy run. ---
$newRows = @(
    @("BYMA.BA", 20, 7, 1000000, 3569549.789, "'2019-08-13", 45190, 24, 10, 14, 0.417, 0.583, -0.075, 369946.83, 0.238, -138062.594, -0.08, -302669.471, 1312411.312, 1.312411312, 3, 4, 61.73913043478261),
    @("PAMP.BA", 20, 7, 1000000, 3457455.821599243, "'2019-08-13", 45190, 32, 14, 18, 0.438, 0.5620000000000001, -0.165, 265916.153, 0.176, -119193.557, -0.08, -241159.874, 1296548.942, 1.296548941523559, 3, 4, 43.03225806451613),
    @("YPFD.BA", 20, 7, 1000000, 5388942.6954, "'2019-08-13", 45190, 27, 14, 13, 0.519, 0.481, -0.232, 435622.156, 0.233, -205632.772, -0.104, -888172.085, 1476139.244, 1.4761392437, 3, 5, 53.69230769230769)
)

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = 3 + $i
    $rowData = $newRows[$i]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $ws.Cells.Item($r, 1 + $c).Value = $rowData[$c]
    }
}

# Copy row 2's per-column formatting (A:W) down into the new rows 3-5, so they
# pick up the same styles (bold/border on A, date format on G, plain text on F)
# instead of minting new style entries (e.g. from the quote-prefixed date text).
$ws.Range("A2:W2").Copy()
$ws.Range("A3:W5").PasteSpecial(-4122)

$wb.Save()
